$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Developer name cell: was an empty bold-styled placeholder cell; typing a
# real name here and the font reverting to the regular (non-bold) weight.
$ws.Range("C3").Value = "Parth Joshi"
$ws.Range("C3").Font.Bold = $false

$ws.Range("E7").Value = 'Valid title, author, and genre values exist.'
$ws.Range("F7").Value = '"One Piece", "Eichiro Oda", Genre.FICTION'
$ws.Range("G7").Value = 'Object created successfully with attributes set to input values.'
$ws.Range("E8").Value = 'Title is blank.'
$ws.Range("F8").Value = '"", "Eichiro Oda", Genre.FICTION'
$ws.Range("G8").Value = 'Raises ValueError with message "Title cannot be blank."'
$ws.Range("E9").Value = 'Author is blank.'
$ws.Range("F9").Value = '"One Piece", "", Genre.FICTION'
$ws.Range("G9").Value = 'Raises ValueError with message "Author cannot be blank."'
$ws.Range("E10").Value = 'Genre is invalid.'
$ws.Range("F10").Value = '"One Piece", "Eichiro Oda", "Invalid"'
$ws.Range("G10").Value = 'Raises ValueError with message "Invalid Genre."'
$ws.Range("E11").Value = 'LibraryItem object created with title set.'
$ws.Range("F11").Value = 'None'
$ws.Range("G11").Value = 'Returns  "One Piece"'
$ws.Range("E12").Value = 'LibraryItem object created with author set.'
$ws.Range("F12").Value = 'None'
$ws.Range("G12").Value = 'Returns ''Eichiro Oda'''
$ws.Range("E13").Value = 'LibraryItem object created with genre set.'
$ws.Range("F13").Value = 'None'
$ws.Range("G13").Value = 'Returns Genre.FICTION'

# Scroll / selection state left behind after data entry.
$ws.Range("J12").Select()

$view = $excel.ActiveWindow
$view.ScrollRow = 7

# Page setup: orientation explicitly set to portrait.
$ws.PageSetup.Orientation = 1
